$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column N (14) width: 11.7109375 -> 12.7109375 (closest achievable via ColumnWidth quantization)
$ws.Columns.Item(14).ColumnWidth = 11.878

# Cell value updates per diff
$ws.Range("A2").Value = 0.75134272579534911
$ws.Range("BL2").Value = 0.95351497333110502
$ws.Range("A3").Value = 0.87623908928610117
$ws.Range("B3").Value = 0.99735318663575923
$ws.Range("AB3").Value = 0.65690987768953346
$ws.Range("F4").Value = 0.97436525659723583
$ws.Range("D5").Value = 0.8791738010364698
$ws.Range("F5").Value = 0.9587440982390063
$ws.Range("BD6").Value = 0.79063186229626536
$ws.Range("E7").Value = 0.94179802443541782
$ws.Range("H7").Value = 0.99212285186500482
$ws.Range("I7").Value = 0.64608166852314497
$ws.Range("F8").Value = 0.83745504698773043
$ws.Range("I8").Value = 0.90840668904710753
$ws.Range("AQ8").Value = 0.87563477761128583
$ws.Range("J9").Value = 0.8929144663700127
$ws.Range("K9").Value = 0.92985592884137702
$ws.Range("AK9").Value = 0.89890338698725025
$ws.Range("AM9").Value = 0.8786365584303959
$ws.Range("K10").Value = 0.98322898721245822
$ws.Range("L10").Value = 0.8726456782468125
$ws.Range("M11").Value = 0.97405841159620143
$ws.Range("BF11").Value = 0.52817448967785596
$ws.Range("BK11").Value = 0.89119971461196656
$ws.Range("K12").Value = 0.79549755564760039
$ws.Range("BP12").Value = 0.78041092495049702
$ws.Range("C13").Value = 0.94192120231323018
$ws.Range("L13").Value = 0.99926830383355747
$ws.Range("N13").Value = 0.76263558616952376
$ws.Range("L14").Value = 0.9876276010007552
$ws.Range("O14").Value = 0.87688635897424816
$ws.Range("BF14").Value = 0.99223254599091892
$ws.Range("P15").Value = 0.74794524313241895
$ws.Range("BB15").Value = 0.98241799885874093
$ws.Range("Y16").Value = 0.67285127217849949
$ws.Range("O17").Value = 0.99653225685572655
$ws.Range("Q18").Value = 0.88379211628472509
$ws.Range("S18").Value = 0.86254057350027269
$ws.Range("T18").Value = 0.92189845089511091
$ws.Range("Q19").Value = 0.84325229091148568
$ws.Range("U19").Value = 0.69994940446090825
$ws.Range("S20").Value = 0.91009086533277794
$ws.Range("AU20").Value = 0.64447697039734519
$ws.Range("V21").Value = 0.88113171465665929
$ws.Range("T22").Value = 0.69308002115992862
$ws.Range("V23").Value = 0.96567733601460581
$ws.Range("V24").Value = 0.83952315329046667
$ws.Range("W24").Value = 0.71784710657428807
$ws.Range("Y24").Value = 0.88937914009878427
$ws.Range("Z24").Value = 0.97198803564070158
$ws.Range("C25").Value = 0.99357348912316124
$ws.Range("AB26").Value = 0.91976760356794629
$ws.Range("AB27").Value = 0.86882371493857669
$ws.Range("AC27").Value = 0.92300515982128228
$ws.Range("Y28").Value = 0.88128876464398065
$ws.Range("AB29").Value = 0.90437379339388368
$ws.Range("AB30").Value = 0.68633263065650207
$ws.Range("AC30").Value = 0.60443473785774016
$ws.Range("AC31").Value = 0.83097306564434748
$ws.Range("AD31").Value = 0.76371742254687547
$ws.Range("AG31").Value = 0.99634481819742682
$ws.Range("AD32").Value = 0.98446965350020577
$ws.Range("AE32").Value = 0.73332531002204671
$ws.Range("AG32").Value = 0.89886825263382941
$ws.Range("AH32").Value = 0.83998381997735849
$ws.Range("K33").Value = 0.80297779699169225
$ws.Range("V33").Value = 0.99869313738396381
$ws.Range("AI34").Value = 0.73789957294357178
$ws.Range("AJ34").Value = 0.65944685574878037
$ws.Range("AK35").Value = 0.99801120204283045
$ws.Range("AU35").Value = 0.88305476532015259
$ws.Range("AL36").Value = 0.85478439409139728
$ws.Range("BM36").Value = 0.91964534645116047
$ws.Range("AJ37").Value = 0.99729873600678642
$ws.Range("AK39").Value = 0.62530325846445833
$ws.Range("AL39").Value = 0.85951365426205273
$ws.Range("AO39").Value = 0.87798348050495112
$ws.Range("I40").Value = 0.89630820531485877
$ws.Range("U40").Value = 0.69639255133073363
$ws.Range("AL40").Value = 0.99327225231563554
$ws.Range("AP41").Value = 0.64761516776219219
$ws.Range("AQ41").Value = 0.97821528185627038
$ws.Range("BA41").Value = 0.70376516325422223
$ws.Range("R42").Value = 0.69377473123242428
$ws.Range("AQ42").Value = 0.9559524056829235
$ws.Range("R43").Value = 0.95886986604845315
$ws.Range("AS44").Value = 0.91594545377228664
$ws.Range("AT44").Value = 0.8481328092336442
$ws.Range("AT45").Value = 0.85034771945454402
$ws.Range("BG46").Value = 0.6428778864699749
$ws.Range("P47").Value = 0.68574619722805785
$ws.Range("V47").Value = 0.98475829979995277
$ws.Range("AS47").Value = 0.75871820096743636
$ws.Range("AT47").Value = 0.82406177843077622
$ws.Range("AW47").Value = 0.94694067614448696
$ws.Range("AY47").Value = 0.80493850536680567
$ws.Range("BH47").Value = 0.85344257222736331
$ws.Range("S48").Value = 0.98501287408416782
$ws.Range("AX48").Value = 0.92091217096662104
$ws.Range("Z49").Value = 0.59788785971010505
$ws.Range("AY49").Value = 0.99168156061351986
$ws.Range("C50").Value = 0.72157240175038706
$ws.Range("AW50").Value = 0.90361694948863613
$ws.Range("AZ50").Value = 0.92161996493949483
$ws.Range("W51").Value = 0.72592021218575353
$ws.Range("AX51").Value = 0.96898682649592671
$ws.Range("AZ51").Value = 0.69868333216052392
$ws.Range("AD52").Value = 0.66657002650343955
$ws.Range("AF52").Value = 0.57717013361779057
$ws.Range("BC53").Value = 0.98178407342304219
$ws.Range("BD54").Value = 0.94687074890705536
$ws.Range("BD55").Value = 0.99416677951953925
$ws.Range("BE56").Value = 0.62856673493998438
$ws.Range("BF56").Value = 0.88084976103840518
$ws.Range("BA57").Value = 0.78703949900303105
$ws.Range("BC57").Value = 0.62816760389009618
$ws.Range("BG57").Value = 0.78284045065203334
$ws.Range("H58").Value = 0.84360954303571389
$ws.Range("BH59").Value = 0.78885891082933868
$ws.Range("BI59").Value = 0.95519411838683088
$ws.Range("AV60").Value = 0.84553742258677622
$ws.Range("BF60").Value = 0.71346669254008366
$ws.Range("BH61").Value = 0.94261196929069491
$ws.Range("BJ61").Value = 0.87710126962281221
$ws.Range("BK62").Value = 0.75213648619413898
$ws.Range("BL62").Value = 0.8942619616115941
$ws.Range("AS63").Value = 0.70191908917155377
$ws.Range("BM64").Value = 0.8719378010891381
$ws.Range("BI66").Value = 0.96974061715647308
$ws.Range("BL66").Value = 0.85977990465878706
$ws.Range("BM66").Value = 0.92072157908165431
$ws.Range("A67").Value = 0.79407302236938859
$ws.Range("BN67").Value = 0.74424428971452294
$ws.Range("AY68").Value = 0.51804225596771136
$ws.Range("BN68").Value = 0.99159377182314534
$ws.Range("BO68").Value = 0.68314318016745912
